# Apply the edits described in the diff to "Hoja1" (sheet1.xml):
#  - A1 date serial 45310 -> 45311 (2024-01-19 -> 2024-01-20)
#  - D29 298.5 -> 300
#  - D30 328 -> 223.526

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45311
$ws.Range("D29").Value = 300
$ws.Range("D30").Value = 223.526
